# Regen sval data to filter save games
# Updates numeric columns B:E and G (sum) for rows 2-8. Column A (dates) and
# column F (Win) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    3 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    4 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 6.48142807727062, 28.30127388105354)
    5 = @(0.00006486019690155054, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 0.9904421852976051)
    6 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
    7 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 3.594575437922795)
    8 = @(1.505614041169197, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 6.741336633845642)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
